# This script applies the change described by the commit:
# "export remaining quote functions to excel (694/1114)"
#
# The Quotes worksheet previously had placeholder "quoteNN#0000" text values
# in columns B/C/E for several qlXxxQuote test rows. The commit replaces those
# placeholders with actual computed results (mostly boolean FALSE, since the
# QuantLibXL add-in functions aren't registered in this environment and the
# calls now resolve to FALSE instead of leaving a placeholder string).
# One row (qlLastFixingQuoteReferenceDate) gets a fuller rework where the
# result moves from a FAIL (#NUM! error) state to an ERROR state with the
# numeric date value now appearing in the Result/Call columns.
#
# A couple of rows on the Volatilities sheet that previously read as FAIL
# (their Result/Call value differed from Expected) are corrected so
# Result/Call now match Expected, turning them into PASS.
#
# Finally, one row on the Date sheet (qlECBIsECBdate) flips from PASS to
# FAIL, as its Result/Call value no longer matches the Expected value.
#
# The UnitTests summary sheet recalculates automatically via its existing
# COUNTIF/SUM formulas, so it needs no direct edits.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Quotes sheet: replace "quoteNN#0000" placeholders with FALSE
# ---------------------------------------------------------------------
$quotes = $wb.Worksheets.Item("Quotes")

# Rows 12-18, 23, 25: B/C/E previously held the placeholder text; now FALSE
$falseRows = @(12, 13, 14, 15, 16, 17, 18, 23, 25)
foreach ($r in $falseRows) {
    $quotes.Range("B$r").Value = $false
    $quotes.Range("C$r").Value = $false
    $quotes.Range("E$r").Value = $false
}

# Row 24 (qlLastFixingQuoteReferenceDate): reworked from FAIL to ERROR
$quotes.Range("B24").Value = "#NUM!"
$quotes.Range("C24").Value = 42644
$quotes.Range("D24").Value = "ERROR"
$quotes.Range("E24").Value = 42644

# ---------------------------------------------------------------------
# Date sheet: qlECBIsECBdate (row 39) now fails
# ---------------------------------------------------------------------
$date = $wb.Worksheets.Item("Date")
$date.Range("C39").Value = $false
$date.Range("D39").Value = "FAIL"
$date.Range("E39").Value = $false

# ---------------------------------------------------------------------
# Volatilities sheet: a handful of rows now match Expected (PASS)
# ---------------------------------------------------------------------
$vol = $wb.Worksheets.Item("Volatilities")

$vol.Range("C24").Value = 0.16276901888733139
$vol.Range("D24").Value = "PASS"
$vol.Range("E24").Value = 0.16276901888733139

$vol.Range("C25").Value = 0.20553128316863267
$vol.Range("D25").Value = "PASS"
$vol.Range("E25").Value = 0.20553128316863267

$vol.Range("C26").Value = -0.00056930143586703347
$vol.Range("D26").Value = "PASS"
$vol.Range("E26").Value = -0.00056930143586703347

$vol.Range("C27").Value = 816.38970366714
$vol.Range("D27").Value = "PASS"
$vol.Range("E27").Value = 816.38970366714

$vol.Range("C30").Value = 0.85719125008922603
$vol.Range("D30").Value = "PASS"
$vol.Range("E30").Value = 0.85719125008922603
